# Apply updated Leve price/profit calculations across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1639.4565
$ws.Range("I40").Value = 1503.3334
$ws.Range("J40").Value = 1832.8948
$ws.Range("K40").Value = 1503.3334
$ws.Range("L40").Value = 1832.8948
$ws.Range("M40").Value = -1328.3334
$ws.Range("N40").Value = -2182.8948

$ws.Range("H74").Value = 5160.533
$ws.Range("I74").Value = 5040
$ws.Range("J74").Value = 5240.8887
$ws.Range("K74").Value = 5040
$ws.Range("L74").Value = 5240.8887
$ws.Range("M74").Value = -4104
$ws.Range("N74").Value = -7112.8887

$ws.Range("H77").Value = 5160.533
$ws.Range("I77").Value = 5040
$ws.Range("J77").Value = 5240.8887
$ws.Range("K77").Value = 25200
$ws.Range("L77").Value = 26204.4435
$ws.Range("M77").Value = -20520
$ws.Range("N77").Value = -35564.4435

$ws.Range("H132").Value = 2595.1853
$ws.Range("I132").Value = 1347.6595
$ws.Range("J132").Value = 10971.429
$ws.Range("K132").Value = 4042.9785
$ws.Range("L132").Value = 32914.287
$ws.Range("M132").Value = -1512.9785
$ws.Range("N132").Value = -37974.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 5500
$ws.Range("I31").Value = 5500
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 5500
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -5206
$ws.Range("N31").Value = ""

$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248

$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96240

$ws.Range("H74").Value = 1665.8572
$ws.Range("I74").Value = 1093.7826
$ws.Range("J74").Value = 3213.8235
$ws.Range("K74").Value = 1093.7826
$ws.Range("L74").Value = 3213.8235
$ws.Range("M74").Value = -219.7826
$ws.Range("N74").Value = -4961.8235

$ws.Range("H77").Value = 1665.8572
$ws.Range("I77").Value = 1093.7826
$ws.Range("J77").Value = 3213.8235
$ws.Range("K77").Value = 5468.913
$ws.Range("L77").Value = 16069.1175
$ws.Range("M77").Value = -1100.913
$ws.Range("N77").Value = -24805.1175

$ws.Range("H110").Value = 2679
$ws.Range("I110").Value = 3866
$ws.Range("J110").Value = 898.5
$ws.Range("K110").Value = 3866
$ws.Range("L110").Value = 898.5
$ws.Range("M110").Value = -1821
$ws.Range("N110").Value = -4988.5

$ws.Range("H132").Value = 19155.295
$ws.Range("I132").Value = 28874.264
$ws.Range("K132").Value = 86622.792
$ws.Range("M132").Value = -84092.792

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1831.8182
$ws.Range("I105").Value = 1765
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1765
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -18
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2301.971
$ws.Range("I31").Value = 1620.0889
$ws.Range("J31").Value = 3580.5
$ws.Range("K31").Value = 1620.0889
$ws.Range("L31").Value = 3580.5
$ws.Range("M31").Value = -1325.0889
$ws.Range("N31").Value = -4170.5

$ws.Range("H34").Value = 2301.971
$ws.Range("I34").Value = 1620.0889
$ws.Range("J34").Value = 3580.5
$ws.Range("K34").Value = 1620.0889
$ws.Range("L34").Value = 3580.5
$ws.Range("M34").Value = -1418.0889
$ws.Range("N34").Value = -3984.5

$ws.Range("H62").Value = 2061362.5
$ws.Range("I62").Value = 7939501.5
$ws.Range("J62").Value = 4013.85
$ws.Range("K62").Value = 7939501.5
$ws.Range("L62").Value = 4013.85
$ws.Range("M62").Value = -7938877.5
$ws.Range("N62").Value = -5261.85

$ws.Range("H65").Value = 2061362.5
$ws.Range("I65").Value = 7939501.5
$ws.Range("J65").Value = 4013.85
$ws.Range("K65").Value = 39697507.5
$ws.Range("L65").Value = 20069.25
$ws.Range("M65").Value = -39694387.5
$ws.Range("N65").Value = -26309.25

$ws.Range("H86").Value = 4093
$ws.Range("I86").Value = 1919.5
$ws.Range("J86").Value = 8440
$ws.Range("K86").Value = 1919.5
$ws.Range("L86").Value = 8440
$ws.Range("M86").Value = -796.5
$ws.Range("N86").Value = -10686

$ws.Range("H89").Value = 4093
$ws.Range("I89").Value = 1919.5
$ws.Range("J89").Value = 8440
$ws.Range("K89").Value = 9597.5
$ws.Range("L89").Value = 42200
$ws.Range("M89").Value = -3981.5
$ws.Range("N89").Value = -53432

$ws.Range("H134").Value = 1716.4
$ws.Range("I134").Value = 1030.762
$ws.Range("K134").Value = 3092.286
$ws.Range("M134").Value = -557.2860000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2456.4443
$ws.Range("I121").Value = 455.8
$ws.Range("J121").Value = 4957.25
$ws.Range("K121").Value = 1367.4
$ws.Range("L121").Value = 14871.75
$ws.Range("M121").Value = -57.40000000000009
$ws.Range("N121").Value = -17491.75

$ws.Range("H129").Value = 1085.375
$ws.Range("I129").Value = 608.3333
$ws.Range("J129").Value = 2516.5
$ws.Range("K129").Value = 1824.9999
$ws.Range("L129").Value = 7549.5
$ws.Range("M129").Value = 3175.0001
$ws.Range("N129").Value = -17549.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3184.7778
$ws.Range("I126").Value = 2746.5833
$ws.Range("J126").Value = 3535.3333
$ws.Range("K126").Value = 8239.749899999999
$ws.Range("L126").Value = 10605.9999
$ws.Range("M126").Value = -5769.749899999999
$ws.Range("N126").Value = -15545.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 954.43475
$ws.Range("I46").Value = 1267.9231
$ws.Range("J46").Value = 546.9
$ws.Range("K46").Value = 1267.9231
$ws.Range("L46").Value = 546.9
$ws.Range("M46").Value = -1079.9231
$ws.Range("N46").Value = -922.9

$ws.Range("H62").Value = 39000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 39000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 39000
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -40248

$ws.Range("H65").Value = 39000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 39000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 117000
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -123240

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4026.8064
$ws.Range("I62").Value = 3605.3572
$ws.Range("J62").Value = 4373.8823
$ws.Range("K62").Value = 3605.3572
$ws.Range("L62").Value = 4373.8823
$ws.Range("M62").Value = -2981.3572
$ws.Range("N62").Value = -5621.8823

$ws.Range("H65").Value = 4026.8064
$ws.Range("I65").Value = 3605.3572
$ws.Range("J65").Value = 4373.8823
$ws.Range("K65").Value = 18026.786
$ws.Range("L65").Value = 21869.4115
$ws.Range("M65").Value = -14906.786
$ws.Range("N65").Value = -28109.4115

